# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F.
$updates = @{
    3  = 3121
    7  = 1728
    9  = 93
    11 = 6
    12 = 1415
    14 = 550
    15 = 353
    16 = 70
    17 = 10
    18 = 79
    21 = 128
    22 = 94
    23 = 112
    24 = 3327
    25 = 401
    26 = 196
    27 = 348
    28 = 24
    31 = 1032
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
